# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated data (gh-pages output at 456a3b4).
#
# Row 2 (丽水·动漫游戏展):           F2  438 -> 439
# Row 3 (丽水·CCAC动漫游戏嘉年华): F3   12 ->  13

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 439
    $ws.Range("F3").Value = 13
}
